$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New competitor analysis content (Simplybook.me strengths/weaknesses, Quandoo strength)
$ws.Range("F9").Value = 'Rather than be a hub for lots of different restaurants, simplybook.me is much more personalized for the specific brand and makes it so the business can integrate it directly into their website or have them create a website for them.'
$ws.Range("F10").Value = 'Can split the restaurant into different sections so that it can be more personalised for the customer. For example,  you can label some tables as "By window", and the user can select tables specific to being next to a window.'
$ws.Range("F11").Value = 'Can set the amount of occupancy you want within the restaurant. This is perfect for when you want to allow walk-ins as well. For example, if you have 40 available seats, you can allow up to 30 up-front bookings for allowance of walk-ins. You also might want to give certain tables more than the average time. Because of this, you dont want the occupancy to be total availability'
$ws.Range("F14").Value = 'Pricing is very expensive in comparison to it''s competitors'
$ws.Range("F15").Value = 'Some consumers say the interface dashboard is quite confusing initially.'
$ws.Range("F16").Value = 'Not specific to restaurant scheduling so doesn''t offer as many CRM capabilities in comparison.'
$ws.Range("G9").Value = 'Customer no-show prevention. Automatic confirmation emails and guest reliability scores.'

# Row heights to fit new wrapped text
$ws.Rows.Item(9).RowHeight = 112.5
$ws.Rows.Item(10).RowHeight = 102.75
$ws.Rows.Item(11).RowHeight = 143.25

# Column F widened to fit new content
$ws.Columns.Item(6).ColumnWidth = 33.59

# Restore active selection
$ws.Range("G9").Select() | Out-Null
